$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) values for rows 2-51
# D column values are prefixed with a leading apostrophe so Excel
# stores them as text (matching the original inline-string format)
# instead of auto-converting them to floating point numbers.
$ws.Range("D2").Value = "'28.661.34"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "'1.875.54"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D5").Value = "'314.66"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'0.5081"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "'0.3914"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "'0.08377"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "'42.20"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D12").Value = "'6.194"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "'1.876.18"
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").Value = "'20.37"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'7.265"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'1.008"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "'93.15"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D19").Value = "'0.06716"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D22").Value = "'5.931"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'28.688.08"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D25").Value = "'2.193"
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").Value = "'2.087.87"
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("D27").Value = "'157.48"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("D28").Value = "'20.62"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "'2.419"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").Value = "'126.56"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "'0.1038"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "'1.047"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D34").Value = "'3.631"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'0.02454"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").Value = "'0.06551"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("D37").Value = "'9.016"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D39").Value = "'5.043"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").Value = "'1.192"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").Value = "'1.241"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").Value = "'0.6389"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D45").Value = "'0.5994"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("D46").Value = "'13.04"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "'3.677"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "'2.004"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D51").Value = "'122.19"
$ws.Range("E51").Value = "  +1.36%  "
$ws.Range("E4").Value = "  -1.19%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  -0.70%  "

# Rows 49 and 50: EOS and WEMIXTOKEN swap positions, with updated price/volume values
$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'1.213"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.221"
$ws.Range("E50").Value = "  +1.62%  "
